# Add the "2022-Q3" worksheet (right after "总计", before "2021-Q4"),
# populate it with the new quarter's holdings data, and insert the
# matching summary row at the top of the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Create the new sheet in the correct position and rename it.
# ---------------------------------------------------------------
$q4sheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($q4sheet, $null)
$newSheet.Name = "2022-Q3"

# re-fetch a fresh handle to "2021-Q4" - the COM layer rebinds the
# variable used for Add()'s "Before" argument to the newly created sheet
$q4sheet = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------
# 2. Bring over the sheet-level look & feel (outline + header style +
#    index-column style + page margins) from the "2021-Q4" sheet so the
#    new tab matches its siblings.
# ---------------------------------------------------------------
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

$q4sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$q4sheet.Range("A2").Copy()
$newSheet.Range("A2:A7").PasteSpecial(-4122)

# ---------------------------------------------------------------
# 3. Header row
# ---------------------------------------------------------------
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

# ---------------------------------------------------------------
# 4. Data rows. Columns B-G are stored as text (fund code / name /
#    scale / position% columns carry trailing zeros etc that must
#    survive as literal text), columns A and H are real numbers.
# ---------------------------------------------------------------
function Set-TextCell($sheet, $row, $col, $text) {
    $cell = $sheet.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# Row 2
$newSheet.Cells.Item(2,1).Value = 0
Set-TextCell $newSheet 2 2 "016935"
Set-TextCell $newSheet 2 3 "景顺长城中证500指数增强C"
Set-TextCell $newSheet 2 4 "15.57"
Set-TextCell $newSheet 2 5 "93.89"
Set-TextCell $newSheet 2 6 "2.06"
Set-TextCell $newSheet 2 7 "0.3207"
$newSheet.Cells.Item(2,8).Value = 5

# Row 3
$newSheet.Cells.Item(3,1).Value = 1
Set-TextCell $newSheet 3 2 "000978"
Set-TextCell $newSheet 3 3 "景顺长城量化精选股票"
Set-TextCell $newSheet 3 4 "7.14"
Set-TextCell $newSheet 3 5 "93.64"
Set-TextCell $newSheet 3 6 "1.92"
Set-TextCell $newSheet 3 7 "0.1371"
$newSheet.Cells.Item(3,8).Value = 7

# Row 4
$newSheet.Cells.Item(4,1).Value = 2
Set-TextCell $newSheet 4 2 "014155"
Set-TextCell $newSheet 4 3 "国泰君安中证500指数增强A"
Set-TextCell $newSheet 4 4 "6.64"
Set-TextCell $newSheet 4 5 "92.15"
Set-TextCell $newSheet 4 6 "1.27"
Set-TextCell $newSheet 4 7 "0.0843"
$newSheet.Cells.Item(4,8).Value = 3

# Row 5
$newSheet.Cells.Item(5,1).Value = 3
Set-TextCell $newSheet 5 2 "014156"
Set-TextCell $newSheet 5 3 "国泰君安中证500指数增强C"
Set-TextCell $newSheet 5 4 "4.02"
Set-TextCell $newSheet 5 5 "92.15"
Set-TextCell $newSheet 5 6 "1.27"
Set-TextCell $newSheet 5 7 "0.0511"
$newSheet.Cells.Item(5,8).Value = 3

# Row 6
$newSheet.Cells.Item(6,1).Value = 4
Set-TextCell $newSheet 6 2 "008851"
Set-TextCell $newSheet 6 3 "景顺长城量化对冲策略三个月定期开放灵活配置混合"
Set-TextCell $newSheet 6 4 "2.96"
Set-TextCell $newSheet 6 5 "64.77"
Set-TextCell $newSheet 6 6 "1.32"
Set-TextCell $newSheet 6 7 "0.0391"
$newSheet.Cells.Item(6,8).Value = 8

# Row 7
$newSheet.Cells.Item(7,1).Value = 5
Set-TextCell $newSheet 7 2 "006682"
Set-TextCell $newSheet 7 3 "景顺长城中证500指数增强A"
Set-TextCell $newSheet 7 4 "0.00"
Set-TextCell $newSheet 7 5 "93.89"
Set-TextCell $newSheet 7 6 "2.06"
$newSheet.Cells.Item(7,7).Value = 0
$newSheet.Cells.Item(7,8).Value = 5

# ---------------------------------------------------------------
# 5. Update the "总计" (summary) sheet: insert a new row right under
#    the header for the 2022-Q3 totals, pushing the older quarters
#    down.
# ---------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q3"
$total.Cells.Item(2,3).Value = 6
$total.Cells.Item(2,4).Value = 0.63

# the row-insert shifted the old index values down along with the rest
# of the row; renumber column A back to a clean 0..n-1 sequence so it
# reads 0,1,2,3 top-to-bottom like the source data.
$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(5,1).Value = 3

# ---------------------------------------------------------------
# 6. Restore the originally-active sheet.
# ---------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
